$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Historias de Usuario")

# --- 1. Give the existing HU-005 row (row 9) the same "wrap text" treatment on
#        the Comentarios cell that the new rows below will use. ---
$ws.Cells.Item(9, 9).WrapText = $true

# --- 2. Fill in five new backlog items directly into the existing template
#        rows 10-14 (they previously held empty "XX-XXXX-XXXX" placeholder
#        rows), then center all 8 columns the way the other filled rows are. ---

# HU-006
$ws.Cells.Item(10, 2).Value = "HU-006"
$ws.Cells.Item(10, 3).Value = "Como asesor estudiantil, quiero recibir una capacitación formal sobre el proceso de baja de materias, para poder informar verbalmente y con precisión a los alumnos."
$ws.Cells.Item(10, 4).Value = "Capacitación a Asesores"
$ws.Cells.Item(10, 5).Value = "Pendiente"
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(10, 7).Value = "Sprint 3"
$ws.Cells.Item(10, 8).Value = "Alta"
$ws.Cells.Item(10, 9).Value = "Asegurar que el 70% de los asesores estén capacitados antes del lanzamiento."

# HU-007
$ws.Cells.Item(11, 2).Value = "HU-007"
$ws.Cells.Item(11, 3).Value = "Como coordinador, quiero medir el porcentaje de estudiantes alcanzados por la campaña, para evaluar el impacto y realizar mejoras."
$ws.Cells.Item(11, 4).Value = "Métricas de Alcance"
$ws.Cells.Item(11, 5).Value = "Pendiente"
$ws.Cells.Item(11, 6).Value = 8
$ws.Cells.Item(11, 7).Value = "Sprint 4"
$ws.Cells.Item(11, 8).Value = "Alta"
$ws.Cells.Item(11, 9).Value = "Se debe contactar al menos al 90% de los estudiantes."

# HU-008
$ws.Cells.Item(12, 2).Value = "HU-008"
$ws.Cells.Item(12, 3).Value = "Como responsable de comunicación, quiero verificar que la difusión cumpla con el reglamento universitario, para prevenir problemas de autorización."
$ws.Cells.Item(12, 4).Value = "Cumplimiento Normativo"
$ws.Cells.Item(12, 5).Value = "Pendiente"
$ws.Cells.Item(12, 6).Value = 5
$ws.Cells.Item(12, 7).Value = "Sprint 2"
$ws.Cells.Item(12, 8).Value = "Alta"
$ws.Cells.Item(12, 9).Value = "Requiere validación formal por parte de la coordinación académica."

# HU-009
$ws.Cells.Item(13, 2).Value = "HU-009"
$ws.Cells.Item(13, 3).Value = "Como diseñador, quiero crear material gráfico claro y atractivo (infografías, carteles), para facilitar la comprensión de fechas y procesos."
$ws.Cells.Item(13, 4).Value = "Material Gráfico"
$ws.Cells.Item(13, 5).Value = "Pendiente"
$ws.Cells.Item(13, 6).Value = 5
$ws.Cells.Item(13, 7).Value = "Sprint 2"
$ws.Cells.Item(13, 8).Value = "Media"
$ws.Cells.Item(13, 9).Value = "Deben aprobarse al menos 3 materiales gráficos."

# HU-010
$ws.Cells.Item(14, 2).Value = "HU-010"
$ws.Cells.Item(14, 3).Value = "Como estudiante, quiero recibir un correo masivo con la información clave, para asegurarme de no perder las fechas importantes."
$ws.Cells.Item(14, 4).Value = "Correo Masivo"
$ws.Cells.Item(14, 5).Value = "Pendiente"
$ws.Cells.Item(14, 6).Value = 8
$ws.Cells.Item(14, 7).Value = "Sprint 3"
$ws.Cells.Item(14, 8).Value = "Alta"
$ws.Cells.Item(14, 9).Value = "Envío a más del 90% de las cuentas institucionales."

$newRows = $ws.Range("B10:I14")
$newRows.HorizontalAlignment = -4108
$ws.Range("B10:I14").RowHeight = 30

# --- 3. Extend the blank template rows (previously 10-18) back out below the
#        new data, by duplicating the still-blank template row 18 six more
#        times (rows 19-24). ---
$template = $ws.Range("B18:I18")
$template.Copy($ws.Range("B19:I19"))
$template.Copy($ws.Range("B20:I20"))
$template.Copy($ws.Range("B21:I21"))
$template.Copy($ws.Range("B22:I22"))
$template.Copy($ws.Range("B23:I23"))
$template.Copy($ws.Range("B24:I24"))

# --- 4. Two extra fully blank (no placeholder text), centred/wrapped rows
#        at the very end of the table. ---
$tail = $ws.Range("B25:I26")
$tail.Value = ""
$tail.Interior.ThemeColor = 0
$tail.VerticalAlignment = -4108
$tail.WrapText = $true

# --- 5. Column H ("Prioridad") now only holds short words (Alta/Media/Baja)
#        so re-fit its width instead of the old fixed custom width. ---
$ws.Columns.Item(8).AutoFit()

# --- 6. Update the view: the user scrolled back up and ended on B18. ---
$ws.Range("B18").Select()

# --- 7. Workbook-level bookkeeping: window size + the local path recorded
#        the last time this file was saved. ---
$wb.Windows.Item(1).Width = 1452
$wb.Windows.Item(1).Height = 786
